# Apply TPM-updated values to the Efnb1-Ephb3 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,7).Value = 6.488422666666668
$ws.Cells.Item(2,8).Value = 19.465268
$ws.Cells.Item(2,9).Value = 0.3444204430827928
$ws.Cells.Item(2,10).Value = 0.361825925629615
$ws.Cells.Item(2,13).Value = 0.1030276666666667
$ws.Cells.Item(2,14).Value = 0.309083
$ws.Cells.Item(2,15).Value = 0.005678588141197309
$ws.Cells.Item(2,16).Value = 0.005869434938871914
$ws.Cells.Item(2,17).Value = 0.6684870476937779
$ws.Cells.Item(2,18).Value = 6.016383429244001
$ws.Cells.Item(2,19).Value = 0.00195582184367587
$ws.Cells.Item(2,20).Value = 0.002123713729680133

# Row 3
$ws.Cells.Item(3,7).Value = 6.488422666666668
$ws.Cells.Item(3,8).Value = 19.465268
$ws.Cells.Item(3,9).Value = 0.3444204430827928
$ws.Cells.Item(3,10).Value = 0.361825925629615
$ws.Cells.Item(3,14).Value = 48.75522599999999
$ws.Cells.Item(3,15).Value = 0.8957491941808339
$ws.Cells.Item(3,16).Value = 0.9258536604633588
$ws.Cells.Item(3,17).Value = 105.4481711656187
$ws.Cells.Item(3,18).Value = 949.0335404905679
$ws.Cells.Item(3,19).Value = 0.3085143343508174
$ws.Cells.Item(3,20).Value = 0.3349978576947221

# Row 4
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,7).Value = 6.488422666666668
$ws.Cells.Item(4,8).Value = 19.465268
$ws.Cells.Item(4,9).Value = 0.3444204430827928
$ws.Cells.Item(4,10).Value = 0.361825925629615
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.7697965
$ws.Cells.Item(4,14).Value = 3.539593
$ws.Cells.Item(4,15).Value = 0.09754608390528599
$ws.Cells.Item(4,16).Value = 0.06721628437535049
$ws.Cells.Item(4,17).Value = 11.48318772598734
$ws.Cells.Item(4,18).Value = 68.89912635592401
$ws.Cells.Item(4,19).Value = 0.03359686543964989
$ws.Cells.Item(4,20).Value = 0.02432059431149462

# Row 5
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,7).Value = 6.488422666666668
$ws.Cells.Item(5,8).Value = 19.465268
$ws.Cells.Item(5,9).Value = 0.3444204430827928
$ws.Cells.Item(5,10).Value = 0.361825925629615
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.01861733333333333
$ws.Cells.Item(5,14).Value = 0.055852
$ws.Cells.Item(5,15).Value = 0.001026133772682911
$ws.Cells.Item(5,16).Value = 0.001060620222418814
$ws.Cells.Item(5,17).Value = 0.1207971275928889
$ws.Cells.Item(5,18).Value = 1.087174148336
$ws.Cells.Item(5,19).Value = 0.000353421448649666
$ws.Cells.Item(5,20).Value = 0.0003837598937181753

# Row 6
$ws.Cells.Item(6,9).Value = 0.4517209651039303
$ws.Cells.Item(6,10).Value = 0.4745489404232121
$ws.Cells.Item(6,13).Value = 0.1030276666666667
$ws.Cells.Item(6,14).Value = 0.309083
$ws.Cells.Item(6,15).Value = 0.005678588141197309
$ws.Cells.Item(6,16).Value = 0.005869434938871914
$ws.Cells.Item(6,17).Value = 0.8767470700661112
$ws.Cells.Item(6,18).Value = 7.890723630595001
$ws.Cells.Item(6,19).Value = 0.002565137315569382
$ws.Cells.Item(6,20).Value = 0.002785334131124648

# Row 7
$ws.Cells.Item(7,9).Value = 0.4517209651039303
$ws.Cells.Item(7,10).Value = 0.4745489404232121
$ws.Cells.Item(7,14).Value = 48.75522599999999
$ws.Cells.Item(7,15).Value = 0.8957491941808339
$ws.Cells.Item(7,16).Value = 0.9258536604633588
$ws.Cells.Item(7,17).Value = 138.2994261926766
$ws.Cells.Item(7,19).Value = 0.4046286904864341
$ws.Cells.Item(7,20).Value = 0.4393628735598393

# Row 8
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,9).Value = 0.4517209651039303
$ws.Cells.Item(8,10).Value = 0.4745489404232121
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.7697965
$ws.Cells.Item(8,14).Value = 3.539593
$ws.Cells.Item(8,15).Value = 0.09754608390528599
$ws.Cells.Item(8,16).Value = 0.06721628437535049
$ws.Cells.Item(8,17).Value = 15.06065260129084
$ws.Cells.Item(8,18).Value = 90.363915607745
$ws.Cells.Item(8,19).Value = 0.04406361116380475
$ws.Cells.Item(8,20).Value = 0.03189741652950788

# Row 9
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,9).Value = 0.4517209651039303
$ws.Cells.Item(9,10).Value = 0.4745489404232121
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.01861733333333333
$ws.Cells.Item(9,14).Value = 0.055852
$ws.Cells.Item(9,15).Value = 0.001026133772682911
$ws.Cells.Item(9,16).Value = 0.001060620222418814
$ws.Cells.Item(9,17).Value = 0.1584301865755556
$ws.Cells.Item(9,18).Value = 1.42587167918
$ws.Cells.Item(9,19).Value = 0.0004635261381220614
$ws.Cells.Item(9,20).Value = 0.0005033162027402795

# Row 10
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.331724
$ws.Cells.Item(10,8).Value = 0.995172
$ws.Cells.Item(10,9).Value = 0.01760867516355742
$ws.Cells.Item(10,10).Value = 0.0184985395557192
$ws.Cells.Item(10,13).Value = 0.1030276666666667
$ws.Cells.Item(10,14).Value = 0.309083
$ws.Cells.Item(10,15).Value = 0.005678588141197309
$ws.Cells.Item(10,16).Value = 0.005869434938871914
$ws.Cells.Item(10,17).Value = 0.03417674969733333
$ws.Cells.Item(10,18).Value = 0.307590747276
$ws.Cells.Item(10,19).Value = 0.00009999241396597275
$ws.Cells.Item(10,20).Value = 0.0001085759743864424

# Row 11
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.331724
$ws.Cells.Item(11,8).Value = 0.995172
$ws.Cells.Item(11,9).Value = 0.01760867516355742
$ws.Cells.Item(11,10).Value = 0.0184985395557192
$ws.Cells.Item(11,14).Value = 48.75522599999999
$ws.Cells.Item(11,15).Value = 0.8957491941808339
$ws.Cells.Item(11,16).Value = 0.9258536604633588
$ws.Cells.Item(11,17).Value = 5.391092863207998
$ws.Cells.Item(11,18).Value = 48.51983576887199
$ws.Cells.Item(11,19).Value = 0.01577295658834862
$ws.Cells.Item(11,20).Value = 0.01712694056088885

# Row 12
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 0.6666666666666666
$ws.Cells.Item(12,7).Value = 0.331724
$ws.Cells.Item(12,8).Value = 0.995172
$ws.Cells.Item(12,9).Value = 0.01760867516355742
$ws.Cells.Item(12,10).Value = 0.0184985395557192
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 1.7697965
$ws.Cells.Item(12,14).Value = 3.539593
$ws.Cells.Item(12,15).Value = 0.09754608390528599
$ws.Cells.Item(12,16).Value = 0.06721628437535049
$ws.Cells.Item(12,17).Value = 0.5870839741659999
$ws.Cells.Item(12,18).Value = 3.522503844996
$ws.Cells.Item(12,19).Value = 0.001717657304965297
$ws.Cells.Item(12,20).Value = 0.001243403095305891

# Row 13
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 0.6666666666666666
$ws.Cells.Item(13,7).Value = 0.331724
$ws.Cells.Item(13,8).Value = 0.995172
$ws.Cells.Item(13,9).Value = 0.01760867516355742
$ws.Cells.Item(13,10).Value = 0.0184985395557192
$ws.Cells.Item(13,11).Value = 1
$ws.Cells.Item(13,12).Value = 0.3333333333333333
$ws.Cells.Item(13,13).Value = 0.01861733333333333
$ws.Cells.Item(13,14).Value = 0.055852
$ws.Cells.Item(13,15).Value = 0.001026133772682911
$ws.Cells.Item(13,16).Value = 0.001060620222418814
$ws.Cells.Item(13,17).Value = 0.006175816282666666
$ws.Cells.Item(13,18).Value = 0.055582346544
$ws.Cells.Item(13,19).Value = 0.00001806885627752904
$ws.Cells.Item(13,20).Value = 0.00001961992513801012

# Row 14
$ws.Cells.Item(14,7).Value = 2.7186785
$ws.Cells.Item(14,8).Value = 5.437357
$ws.Cells.Item(14,9).Value = 0.1443137264130649
$ws.Cells.Item(14,10).Value = 0.1010711349827635
$ws.Cells.Item(14,13).Value = 0.1030276666666667
$ws.Cells.Item(14,14).Value = 0.309083
$ws.Cells.Item(14,15).Value = 0.005678588141197309
$ws.Cells.Item(14,16).Value = 0.005869434938871914
$ws.Cells.Item(14,17).Value = 0.2800991022718334
$ws.Cells.Item(14,18).Value = 1.680594613631
$ws.Cells.Item(14,19).Value = 0.0008194982154212233
$ws.Cells.Item(14,20).Value = 0.0005932304509792712

# Row 15
$ws.Cells.Item(15,7).Value = 2.7186785
$ws.Cells.Item(15,8).Value = 5.437357
$ws.Cells.Item(15,9).Value = 0.1443137264130649
$ws.Cells.Item(15,10).Value = 0.1010711349827635
$ws.Cells.Item(15,14).Value = 48.75522599999999
$ws.Cells.Item(15,15).Value = 0.8957491941808339
$ws.Cells.Item(15,16).Value = 0.9258536604633588
$ws.Cells.Item(15,17).Value = 44.183261562947
$ws.Cells.Item(15,18).Value = 265.099569377682
$ws.Cells.Item(15,19).Value = 0.1292689041437362
$ws.Cells.Item(15,20).Value = 0.09357708029097779

# Row 16
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,7).Value = 2.7186785
$ws.Cells.Item(16,8).Value = 5.437357
$ws.Cells.Item(16,9).Value = 0.1443137264130649
$ws.Cells.Item(16,10).Value = 0.1010711349827635
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 1.7697965
$ws.Cells.Item(16,14).Value = 3.539593
$ws.Cells.Item(16,15).Value = 0.09754608390528599
$ws.Cells.Item(16,16).Value = 0.06721628437535049
$ws.Cells.Item(16,17).Value = 4.81150769392525
$ws.Cells.Item(16,18).Value = 19.246030775701
$ws.Cells.Item(16,19).Value = 0.01407723886537332
$ws.Cells.Item(16,20).Value = 0.006793626151140863

# Row 17
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,7).Value = 2.7186785
$ws.Cells.Item(17,8).Value = 5.437357
$ws.Cells.Item(17,9).Value = 0.1443137264130649
$ws.Cells.Item(17,10).Value = 0.1010711349827635
$ws.Cells.Item(17,11).Value = 1
$ws.Cells.Item(17,12).Value = 0.3333333333333333
$ws.Cells.Item(17,13).Value = 0.01861733333333333
$ws.Cells.Item(17,14).Value = 0.055852
$ws.Cells.Item(17,15).Value = 0.001026133772682911
$ws.Cells.Item(17,16).Value = 0.001060620222418814
$ws.Cells.Item(17,17).Value = 0.05061454386066667
$ws.Cells.Item(17,18).Value = 0.303687263164
$ws.Cells.Item(17,19).Value = 0.0001480851885341677
$ws.Cells.Item(17,20).Value = 0.0001071980896655405

# Row 18
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 0.790022
$ws.Cells.Item(18,8).Value = 2.370066
$ws.Cells.Item(18,9).Value = 0.04193619023665445
$ws.Cells.Item(18,10).Value = 0.04405545940869034
$ws.Cells.Item(18,13).Value = 0.1030276666666667
$ws.Cells.Item(18,14).Value = 0.309083
$ws.Cells.Item(18,15).Value = 0.005678588141197309
$ws.Cells.Item(18,16).Value = 0.005869434938871914
$ws.Cells.Item(18,17).Value = 0.08139412327533334
$ws.Cells.Item(18,18).Value = 0.732547109478
$ws.Cells.Item(18,19).Value = 0.0002381383525648604
$ws.Cells.Item(18,20).Value = 0.0002585806527014205

# Row 19
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 0.790022
$ws.Cells.Item(19,8).Value = 2.370066
$ws.Cells.Item(19,9).Value = 0.04193619023665445
$ws.Cells.Item(19,10).Value = 0.04405545940869034
$ws.Cells.Item(19,14).Value = 48.75522599999999
$ws.Cells.Item(19,15).Value = 0.8957491941808339
$ws.Cells.Item(19,16).Value = 0.9258536604633588
$ws.Cells.Item(19,17).Value = 12.839233718324
$ws.Cells.Item(19,18).Value = 115.553103464916
$ws.Cells.Item(19,19).Value = 0.03756430861149738
$ws.Cells.Item(19,20).Value = 0.04078890835693087

# Row 20
$ws.Cells.Item(20,4).Value = "MuSCs"
$ws.Cells.Item(20,5).Value = 3
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 0.790022
$ws.Cells.Item(20,8).Value = 2.370066
$ws.Cells.Item(20,9).Value = 0.04193619023665445
$ws.Cells.Item(20,10).Value = 0.04405545940869034
$ws.Cells.Item(20,11).Value = 2
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = 1.7697965
$ws.Cells.Item(20,14).Value = 3.539593
$ws.Cells.Item(20,15).Value = 0.09754608390528599
$ws.Cells.Item(20,16).Value = 0.06721628437535049
$ws.Cells.Item(20,17).Value = 1.398178170523
$ws.Cells.Item(20,18).Value = 8.389069023138
$ws.Cells.Item(20,19).Value = 0.00409071113149273
$ws.Cells.Item(20,20).Value = 0.00296124428790124

# Row 21
$ws.Cells.Item(21,4).Value = "Resolving-Mac"
$ws.Cells.Item(21,5).Value = 3
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = 0.790022
$ws.Cells.Item(21,8).Value = 2.370066
$ws.Cells.Item(21,9).Value = 0.04193619023665445
$ws.Cells.Item(21,10).Value = 0.04405545940869034
$ws.Cells.Item(21,11).Value = 1
$ws.Cells.Item(21,12).Value = 0.3333333333333333
$ws.Cells.Item(21,13).Value = 0.01861733333333333
$ws.Cells.Item(21,14).Value = 0.055852
$ws.Cells.Item(21,15).Value = 0.001026133772682911
$ws.Cells.Item(21,16).Value = 0.001060620222418814
$ws.Cells.Item(21,17).Value = 0.01470810291466667
$ws.Cells.Item(21,18).Value = 0.132372926232
$ws.Cells.Item(21,19).Value = 0.00004303214109948648
$ws.Cells.Item(21,20).Value = 0.00004672611115680816

